# Workbook was "finalized": the 10-years daily data table (rows 9:39 of the
# original "Data Harian - Table" sheet) is duplicated onto a brand-new
# "Sheet1" worksheet (re-based to rows 1:31), which becomes the new active /
# selected tab. The original sheet keeps all of its data untouched but the
# on-screen selection moves to the data table and it is no longer the
# selected tab.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Keep gridlines visible on the original sheet (matches the source file).
$ws1.Activate()
$excel.ActiveWindow.DisplayGridlines = $true

# Add the new worksheet right after the existing one.
$newSheet = $wb.Worksheets.Add($null, $ws1)

# Copy the full data table (header row 9 + 30 daily rows 10:39) including
# formatting/styles onto the new sheet, anchored at A1.
$src = $ws1.Range("A9:K39")
$src.Copy($newSheet.Range("A1"))

# The data rows wrap text across two lines in real Excel once re-laid-out;
# approximate that taller row height for the 30 data rows (header stays at
# the default height).
$newSheet.Range("A2:K31").RowHeight = 28.8

# Restore the selection on the original sheet to the data table, then make
# the new sheet the active / selected tab (mirrors the authored workbook).
$ws1.Range("A9:K39").Select() | Out-Null
$newSheet.Range("A1:K31").Select() | Out-Null
$newSheet.Activate()
